# Daily attendance processing - 2026-01-01 17:05:12
# Reorders the "Recorded By" (column G) contributor list for every row whose
# list contains a "System" entry alongside other contributors: the list is
# reversed (so "System" moves from one end of the list to the other).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*,*") {
        $parts = $val -split ", "

        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.ToLower() -eq "system") {
                $hasSystem = $true
            }
        }

        if ($hasSystem) {
            $n = $parts.Count
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $newVal = $reversed -join ", "
            $cell.Value = $newVal
            $changed = $changed + 1
        }
    }
}

Write-Output "Rows updated:"
Write-Output $changed
